$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 70

$ws.Cells.Item($row, 1).Value = "T2DX9Z"
$ws.Cells.Item($row, 2).Value = "Goma de rodillo de papel para Samsung"
$ws.Cells.Item($row, 3).Value = "CLX 2160 3160, CLP 300 350, ML 1610 1615 1640 1641 2010 2015 2240 2241 2245, SCX 4321 4521, XP 3117 3122 3124 6110, X WC PE220"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 50000
$ws.Cells.Item($row, 6).Value = 8
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E70-D70)*G70"
$ws.Cells.Item($row, 9).Formula = "=D70*F70"
$ws.Cells.Item($row, 10).Value = 0
